# Rename the left-hand ("_old") and right-hand ("_new") header-suffix
# groups to the respective format-version names ("_FV2404" / "_FV2410"),
# matching the new input file naming scheme.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = $cell.Value2.Replace("_old", "_FV2404")
}

for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = $cell.Value2.Replace("_new", "_FV2410")
}

# Turn the data range into an Excel table ("Table1") covering the full
# used range, header row included.
$rng = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row so it stays visible while scrolling.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
